$d = $word.ActiveDocument

# 1) Fix typo: "avverkningsamnälda" -> "avverkningsanmälda" in the main body text.
$d.Content.Find.Execute(
    "avverkningsamnälda", $true, $false, $false, $false, $false,
    $true, 1, $false, "avverkningsanmälda", 2) | Out-Null

# 2) Update the date in the first-page header (header3.xml) from 2023-11-03 to 2023-11-13,
#    touching only the date text itself (leave the leading tabs/breaks alone).
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage
$hdrFound = $hdr.Range.Find.Execute(
    "2023-11-03", $false, $false, $false, $false, $false,
    $true, 1, $false, "2023-11-13", 2)
if (-not $hdrFound) {
    # Fallback: scan every header/footer story range in case the date text
    # is not where we expect it.
    foreach ($story in $d.StoryRanges) {
        $story.Find.Execute(
            "2023-11-03", $false, $false, $false, $false, $false,
            $true, 1, $false, "2023-11-13", 2) | Out-Null
    }
}

# 3) Document-wide default language switch from en-US to sv-SE (eastAsia/bidi unchanged).
#    Apply to every style's run properties (skips list/numbering-type styles, which
#    cannot carry character formatting).
foreach ($s in $d.Styles) {
    if ($s.Type -ne 4) {
        $fnt = $s.Font
        $fnt.LanguageID = "sv-SE"
        $fnt.LanguageIDFarEast = "en-US"
        $fnt.LanguageIDOther = "ar-SA"
    }
}
